# This script re-orders the data rows 34-50 of the "Artfynd" worksheet so that
# each row's content is replaced by the content that (before this edit) lived
# in a different row, per the permutation observed in the source diff. Row 46
# keeps its own data (it is not touched by the diff).
#
# Only the columns whose values actually differ between rows in this block are
# touched (A, B, D, E, F, G, H, Q, R, Z, AB, AW, AX). Every other column in the
# range (C, I, M, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AY) holds the exact
# same value on every one of these rows, so permuting the rows would be a
# no-op for them - they are deliberately left untouched to avoid Excel's
# automatic type coercion (e.g. turning the literal text "2023-08-22" into a
# real date) or the blanking-out of already-empty cells that can happen when
# round-tripping values through COM unnecessarily.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 34
$lastRow  = 50

# Mapping of sheet row -> sheet row whose original content should end up there.
$rowMap = @{
    34 = 40
    35 = 36
    36 = 49
    37 = 42
    38 = 48
    39 = 43
    40 = 35
    41 = 44
    42 = 50
    43 = 47
    44 = 37
    45 = 39
    46 = 46
    47 = 45
    48 = 34
    49 = 41
    50 = 38
}

$columns = @("A","B","D","E","F","G","H","Q","R","Z","AB","AW","AX")

foreach ($col in $columns) {
    $addr = "$col$firstRow`:$col$lastRow"
    $range = $ws.Range($addr)

    # Excel COM returns a 2-D SAFEARRAY that is 1-based: dim 1 = rows, dim 2 = columns.
    $src = $range.Value2
    $rowCount = $src.GetLength(0)

    # .NET arrays created with New-Object are 0-based, unlike the SAFEARRAY above.
    $dst = New-Object 'object[,]' $rowCount,1

    for ($sheetRow = $firstRow; $sheetRow -le $lastRow; $sheetRow++) {
        $destArrRow = $sheetRow - $firstRow          # 0-based index into $dst
        $sourceRow  = $rowMap[$sheetRow]
        $srcArrRow  = $sourceRow - $firstRow + 1     # 1-based index into $src

        $dst[$destArrRow, 0] = $src[$srcArrRow, 1]
    }

    $range.Value2 = $dst
}
